$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.987.49'
$ws.Range("E2").Value = '  +0.20%  '
$ws.Range("D3").Value = '1.633.04'
$ws.Range("E3").Value = '  -0.63%  '
$c = $ws.Range("D4")
$c.Value = "'0.997"
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.35%  '
$c = $ws.Range("D5")
$c.Value = "'212.01"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("E6").Value = '  -0.54%  '
$c = $ws.Range("D7")
$c.Value = "'0.996"
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.38%  '
$c = $ws.Range("D8")
$c.Value = "'23.51"
$c.Style = "Normal"
$ws.Range("E8").Value = '  -0.76%  '
$ws.Range("E9").Value = '  -2.22%  '
$ws.Range("E10").Value = '  -0.41%  '
$c = $ws.Range("D11")
$c.Value = "'0.0880"
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("D12").Value = '1.862.88'
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("D13").Value = '1.628.94'
$ws.Range("E13").Value = '  -0.86%  '
$ws.Range("E14").Value = '  -0.30%  '
$ws.Range("E15").Value = '  -1.51%  '
$c = $ws.Range("D16")
$c.Value = "'65.66"
$c.Style = "Normal"
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = '27.972.54'
$ws.Range("E17").Value = '  +0.22%  '
$c = $ws.Range("D18")
$c.Value = "'232.18"
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.77%  '
$ws.Range("E19").Value = '  +0.27%  '
$c = $ws.Range("D20")
$c.Value = "'7.56"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.95%  '
$ws.Range("E21").Value = '  -0.43%  '
$c = $ws.Range("D22")
$c.Value = "'10.50"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -4.03%  '
$ws.Range("E23").Value = '  -0.85%  '
$ws.Range("E24").Value = '  -3.59%  '
$c = $ws.Range("D25")
$c.Value = "'154.33"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +1.37%  '
$c = $ws.Range("D26")
$c.Value = "'6.95"
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.22%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range("D27")
$c.Value = "'15.67"
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.44%  '
$ws.Range("B28").Value = 'Stellar'
$ws.Range("C28").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D28")
$c.Value = "'0.111"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.72%  '
$c = $ws.Range("D29")
$c.Value = "'0.997"
$c.Style = "Normal"
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("E30").Value = '  -0.48%  '
$c = $ws.Range("D31")
$c.Value = "'0.0482"
$c.Style = "Normal"
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("E32").Value = '  +1.86%  '
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("D34").Value = '1.408.95'
$ws.Range("E34").Value = '  -1.27%  '
$c = $ws.Range("D36")
$c.Value = "'1.01"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +8.71%  '
$ws.Range("E37").Value = '  +0.46%  '
$ws.Range("E38").Value = '  +1.79%  '
$c = $ws.Range("D39")
$c.Value = "'0.557"
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.44%  '
$ws.Range("E40").Value = '  -2.52%  '
$ws.Range("E41").Value = '  -1.19%  '
$ws.Range("E42").Value = '  -0.34%  '
$c = $ws.Range("D43")
$c.Value = "'66.91"
$c.Style = "Normal"
$ws.Range("E43").Value = '  -2.40%  '
$ws.Range("E44").Value = '  +0.51%  '
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("E46").Value = '  -0.50%  '
$ws.Range("D47").Value = '1.773.70'
$ws.Range("E47").Value = '  -0.61%  '
$c = $ws.Range("D48")
$c.Value = "'88.20"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.81%  '
$ws.Range("E49").Value = '  -3.51%  '
$c = $ws.Range("D50")
$c.Value = "'0.1000"
$c.Style = "Normal"
$ws.Range("E50").Value = '  -1.00%  '
$c = $ws.Range("D51")
$c.Value = "'0.0504"
$c.Style = "Normal"
$ws.Range("E51").Value = '  -0.36%  '
